# Update employee absence data rows 2-11 with new values per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 27665
$ws.Cells.Item(2, 2).Value = "Aylla Ribeiro"
$ws.Cells.Item(2, 3).Value = "Operacoes"
$ws.Cells.Item(2, 4).Value = "Consulta medica"
$ws.Cells.Item(2, 5).Value = 7
$ws.Cells.Item(2, 6).Value = 45078
$ws.Cells.Item(2, 7).Value = 6609.42

# Row 3
$ws.Cells.Item(3, 1).Value = 86218
$ws.Cells.Item(3, 2).Value = "Leonardo Viana"
$ws.Cells.Item(3, 3).Value = "Engenharia"
$ws.Cells.Item(3, 4).Value = "Problemas pessoais"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 45105
$ws.Cells.Item(3, 7).Value = 3954.66

# Row 4
$ws.Cells.Item(4, 1).Value = 98999
$ws.Cells.Item(4, 2).Value = "Isabelly Fernandes"
$ws.Cells.Item(4, 3).Value = "P&D"
$ws.Cells.Item(4, 4).Value = "Consulta medica"
$ws.Cells.Item(4, 5).Value = 6
$ws.Cells.Item(4, 6).Value = 45081
$ws.Cells.Item(4, 7).Value = 3145.08

# Row 5
$ws.Cells.Item(5, 1).Value = 9042
$ws.Cells.Item(5, 2).Value = "Cecília Silva"
$ws.Cells.Item(5, 3).Value = "Operacoes"
$ws.Cells.Item(5, 4).Value = "Consulta medica"
$ws.Cells.Item(5, 5).Value = 7
$ws.Cells.Item(5, 6).Value = 45091
$ws.Cells.Item(5, 7).Value = 9650.959999999999

# Row 6
$ws.Cells.Item(6, 1).Value = 31093
$ws.Cells.Item(6, 2).Value = "Dr. Marcos Vinicius Ramos"
$ws.Cells.Item(6, 3).Value = "Vendas"
$ws.Cells.Item(6, 4).Value = "Consulta medica"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 45085
$ws.Cells.Item(6, 7).Value = 9436.040000000001

# Row 7
$ws.Cells.Item(7, 1).Value = 89341
$ws.Cells.Item(7, 2).Value = "Arthur Gabriel Campos"
$ws.Cells.Item(7, 3).Value = "Recursos Humanos"
$ws.Cells.Item(7, 4).Value = "Consulta medica"
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 45099
$ws.Cells.Item(7, 7).Value = 4233.34

# Row 8
$ws.Cells.Item(8, 1).Value = 29693
$ws.Cells.Item(8, 2).Value = "Rodrigo da Mata"
$ws.Cells.Item(8, 3).Value = "Financeiro"
$ws.Cells.Item(8, 4).Value = "Outros"
$ws.Cells.Item(8, 5).Value = 7
$ws.Cells.Item(8, 6).Value = 45083
$ws.Cells.Item(8, 7).Value = 6309.7

# Row 9
$ws.Cells.Item(9, 1).Value = 21651
$ws.Cells.Item(9, 2).Value = "André Ribeiro"
$ws.Cells.Item(9, 3).Value = "Marketing"
$ws.Cells.Item(9, 4).Value = "Problemas pessoais"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 45088
$ws.Cells.Item(9, 7).Value = 7407.78

# Row 10
$ws.Cells.Item(10, 1).Value = 40132
$ws.Cells.Item(10, 2).Value = "Srta. Ayla Dias"
$ws.Cells.Item(10, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(10, 4).Value = "Problemas pessoais"
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = 45081
$ws.Cells.Item(10, 7).Value = 8216.52

# Row 11
$ws.Cells.Item(11, 1).Value = 9813
$ws.Cells.Item(11, 2).Value = "Júlia Silva"
$ws.Cells.Item(11, 3).Value = "Financeiro"
$ws.Cells.Item(11, 4).Value = "Doenca"
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 45098
$ws.Cells.Item(11, 7).Value = 4322.96
